# Updated message sequence diagram. Added CONVEYOR_OCCUPIED and
# CONVEYOR_FREE to the messages spreadsheet (Sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 17: CONVEYOR_OCCUPIED message
$ws.Range("A17").Value = "CONVEYOR_OCCUPIED"
$ws.Range("B17").Value = 16
$ws.Range("C17").Value = "Conveyor"
$ws.Range("D17").Value = "N/A"
$ws.Range("E17").Value = "N/A"
$ws.Range("F17").Value = "Sent when the conveyor is occupied"

# New row 18: CONVEYOR_FREE message
$ws.Range("A18").Value = "CONVEYOR_FREE"
$ws.Range("B18").Value = 17
$ws.Range("C18").Value = "Conveyor"
$ws.Range("D18").Value = "N/A"
$ws.Range("E18").Value = "N/A"
$ws.Range("F18").Value = "Sent when the conveyor is free."

# Move the active selection to the last edited cell, as in the authored edit
$ws.Range("F18").Select()
